# Fix wrong title in day9
#
# The title slide's headline still described an older lecture
# ("Bias-Variance Tradeoff, Clustering, K-Means") even though this
# deck is actually about Unsupervised Learning / Clustering / K-Means
# (see slide 2's outline: "Unsupervised Learning and Clustering" /
# "Intro to K-Means"). Correct the headline to read
# "Unsupervised Learning, Clustering, K-Means".
#
# While we're touching that text box we also coalesce the adjacent
# "Sign-In: [See board]" runs on the same slide the same way
# PowerPoint's editor does when text is retyped (purely a run-split
# cleanup - the rendered text is unchanged).

$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)

# --- Title headline fix -------------------------------------------------
$title = $s1.Shapes.Item("Title 1")
$titleRange = $title.TextFrame.TextRange
$titleText = $titleRange.Text

$oldRun1 = "Bias-"
$oldRun2 = "Variance Tradeoff, C"
$oldRun3 = "lustering, K-Means"

$idxRun1 = $titleText.IndexOf($oldRun1)
$idxRun2 = $titleText.IndexOf($oldRun2)
$idxRun3 = $titleText.IndexOf($oldRun3)

if ($idxRun1 -ge 0 -and $idxRun2 -ge 0 -and $idxRun3 -ge 0) {
    # Work right-to-left so earlier offsets stay valid as lengths change.
    $run3 = $titleRange.Characters($idxRun3 + 1, $oldRun3.Length)
    $run3.Text = "Clustering, K-Means"

    $run2 = $titleRange.Characters($idxRun2 + 1, $oldRun2.Length)
    $run2.Text = ""

    $run1 = $titleRange.Characters($idxRun1 + 1, $oldRun1.Length)
    $run1.Text = "Unsupervised Learning, "
}

# --- Sign-In run cleanup -------------------------------------------------
$subtitle = $s1.Shapes.Item("Subtitle 2")
$subRange = $subtitle.TextFrame.TextRange
$subText = $subRange.Text

$signInIdx = $subText.IndexOf("Sign-In")
if ($signInIdx -ge 0) {
    $afterSignIn = $signInIdx + "Sign-In".Length

    $oldColon = ": "
    $oldSee = "[See "
    $oldBoard = "board]"

    $idxColon = $subText.IndexOf($oldColon, $afterSignIn)
    $idxSee = $subText.IndexOf($oldSee, $afterSignIn)
    $idxBoard = $subText.IndexOf($oldBoard, $afterSignIn)

    if ($idxColon -ge 0 -and $idxSee -ge 0 -and $idxBoard -ge 0) {
        $boardRun = $subRange.Characters($idxBoard + 1, $oldBoard.Length)
        $boardRun.Text = ""

        $seeRun = $subRange.Characters($idxSee + 1, $oldSee.Length)
        $seeRun.Text = ""

        $colonRun = $subRange.Characters($idxColon + 1, $oldColon.Length)
        $colonRun.Text = ": [See board]"
    }
}
